$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@(2, "Bitcoin", "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc", "30.650.21", "  +0.76%  ")
    ,@(3, "Ethereum", "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth", "1.880.35", "  -0.35%  ")
    ,@(4, "TetherUSD", "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt", "1.001", "  +0.13%  ")
    ,@(5, "BNB", "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb", "239.24", "  -0.03%  ")
    ,@(6, "USDC", "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc", "1.001", "  +0.16%  ")
    ,@(7, "XRP", "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp", "0.4827", "  -0.08%  ")
    ,@(8, "LidoStakedEther", "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth", "1.880.38", "  -0.17%  ")
    ,@(9, "Cardano", "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada", "0.2831", "  -1.58%  ")
    ,@(10, "Dogecoin", "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge", "0.06533", "  -1.23%  ")
    ,@(11, "WrappedEther", "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth", "1.955.38", "  +3.67%  ")
    ,@(12, "TRON", "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx", "0.07523", "  +1.48%  ")
    ,@(13, "Solana", "https://coinranking.com/coin/zNZHO_Sjf+solana-sol", "16.58", "  -2.20%  ")
    ,@(14, "Polkadot", "https://coinranking.com/coin/25W7FG7om+polkadot-dot", "5.104", "  -1.40%  ")
    ,@(15, "Litecoin", "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc", "88.89", "  +0.07%  ")
    ,@(16, "Polygon", "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic", "0.6662", "  +0.50%  ")
    ,@(17, "WrappedBTC", "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc", "30.608.73", "  +0.78%  ")
    ,@(18, "WrappedliquidstakedEther2.0", "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth", "2.279.58", "  +6.84%  ")
    ,@(19, "Avalanche", "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax", "13.33", "  -1.44%  ")
    ,@(20, "Dai", "https://coinranking.com/coin/MoTuySvg7+dai-dai", "1.001", "  +0.17%  ")
    ,@(21, "ShibaInu", "https://coinranking.com/coin/xz24e0BjL+shibainu-shib", "0.000007619", "  -2.01%  ")
    ,@(22, "BitcoinCash", "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch", "231.25", "  +10.69%  ")
    ,@(23, "Uniswap", "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni", "5.292", "  -2.62%  ")
    ,@(24, "BinanceUSD", "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd", "1.001", "  +0.23%  ")
    ,@(25, "Chainlink", "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link", "6.186", "  -0.17%  ")
    ,@(26, "Cosmos", "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom", "9.339", "  -1.04%  ")
    ,@(27, "Monero", "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr", "167.67", "  +1.64%  ")
    ,@(28, "EthereumClassic", "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc", "18.64", "  +1.12%  ")
    ,@(29, "LidoDAOToken", "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo", "1.946", "  +0.29%  ")
    ,@(30, "Toncoin", "https://coinranking.com/coin/67YlI0K1b+toncoin-ton", "1.425", "  -0.81%  ")
    ,@(31, "Stellar", "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm", "0.09553", "  +4.11%  ")
    ,@(32, "InternetComputer(DFINITY)", "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp", "4.349", "  +1.16%  ")
    ,@(33, "Filecoin", "https://coinranking.com/coin/ymQub4fuB+filecoin-fil", "4.043", "  -0.35%  ")
    ,@(34, "Hedera", "https://coinranking.com/coin/jad286TjB+hedera-hbar", "0.05024", "  -0.84%  ")
    ,@(35, "ARBITRUM", "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb", "1.214", "  +5.10%  ")
    ,@(36, "ImmutableX", "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx", "0.7485", "  +0.33%  ")
    ,@(37, "HuobiToken", "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht", "2.702", "  -0.10%  ")
    ,@(38, "VeChain", "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet", "0.01853", "  -2.81%  ")
    ,@(39, "MXToken", "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx", "2.627", "  -0.47%  ")
    ,@(40, "RenderToken", "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr", "2.092", "  +0.55%  ")
    ,@(41, "TrustWalletToken", "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt", "0.9172", "  -0.01%  ")
    ,@(42, "Quant", "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt", "106.36", "  -0.89%  ")
    ,@(43, "TheSandbox", "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand", "0.4286", "  -1.34%  ")
    ,@(44, "FraxShare", "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs", "5.825", "  -4.18%  ")
    ,@(45, "PaxDollar", "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp", "1.000", "  -0.02%  ")
    ,@(46, "Aptos", "https://coinranking.com/coin/HGYj5JCv5+aptos-apt", "7.427", "  -2.90%  ")
    ,@(47, "Aave", "https://coinranking.com/coin/ixgUfzmLR+aave-aave", "64.75", "  -1.18%  ")
    ,@(48, "Algorand", "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo", "0.1287", "  -4.48%  ")
    ,@(49, "NEARProtocol", "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near", "1.480", "  -7.68%  ")
    ,@(50, "EnergySwap", "https://coinranking.com/coin/SbWqqTui-+energyswap-ens", "8.956", "  +0.17%  ")
    ,@(51, "Elrond", "https://coinranking.com/coin/omwkOTglq+elrond-egld", "33.96", "  -1.08%  ")
)

foreach ($item in $data) {
    $r = $item[0]
    $ws.Cells.Item($r, 2).Value = $item[1]
    $ws.Cells.Item($r, 3).Value = $item[2]
    $ws.Cells.Item($r, 4).Value = $item[3]
    $ws.Cells.Item($r, 5).Value = $item[4]
}

Write-Host "Done updating cryptos list"